$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "40.764.15"
$ws.Range("E2").Value = "  -6.86%  "

# Row 3
$ws.Range("D3").Value = "2.184.79"
$ws.Range("E3").Value = "  -7.12%  "

# Row 4
$ws.Range("E4").Value = "  -0.29%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.27"
$ws.Range("E5").Value = "  -0.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.619"
$ws.Range("E6").Value = "  -7.32%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "69.48"
$ws.Range("E7").Value = "  -4.97%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  -10.76%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.75"
$ws.Range("E10").Value = "  +4.07%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("E11").Value = "  -7.82%  "

# Row 12
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "57.48"
$ws.Range("E12").Value = "  -5.42%  "

# Row 13
$ws.Range("E13").Value = "  -4.43%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.53"
$ws.Range("E14").Value = "  -9.14%  "

# Row 15
$ws.Range("D15").Value = "2.507.11"
$ws.Range("E15").Value = "  -7.37%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.57"
$ws.Range("E16").Value = "  -10.00%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.831"
$ws.Range("E17").Value = "  -8.84%  "

# Row 18
$ws.Range("D18").Value = "2.175.29"
$ws.Range("E18").Value = "  -7.86%  "

# Row 19
$ws.Range("D19").Value = "40.688.33"
$ws.Range("E19").Value = "  -7.01%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  -8.86%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.16"
$ws.Range("E21").Value = "  -7.10%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.04"
$ws.Range("E22").Value = "  -7.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "229.79"
$ws.Range("E23").Value = "  -8.97%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.00"
$ws.Range("E24").Value = "  +6.73%  "

# Row 25
$ws.Range("E25").Value = "  +0.11%  "

# Row 26
$ws.Range("E26").Value = "  -4.96%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.38"
$ws.Range("E27").Value = "  -4.32%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.18"
$ws.Range("E28").Value = "  -5.09%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.66"
$ws.Range("E29").Value = "  -7.73%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.47"
$ws.Range("E30").Value = "  -4.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.16"
$ws.Range("E31").Value = "  -9.34%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.117"
$ws.Range("E32").Value = "  -8.96%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.123"
$ws.Range("E33").Value = "  -7.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0698"
$ws.Range("E34").Value = "  -6.32%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.11"
$ws.Range("E35").Value = "  -4.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.56"
$ws.Range("E36").Value = "  -9.46%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.79"
$ws.Range("E37").Value = "  +0.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.21"
$ws.Range("E38").Value = "  +17.96%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("E39").Value = "  -6.62%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0265"
$ws.Range("E40").Value = "  -4.21%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.80"
$ws.Range("E41").Value = "  -11.89%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "63.79"
$ws.Range("E42").Value = "  -1.87%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.78"
$ws.Range("E43").Value = "  -13.28%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.58"
$ws.Range("E44").Value = "  -4.88%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.190"
$ws.Range("E45").Value = "  -5.76%  "

# Row 46
$ws.Range("E46").Value = "  -0.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0982"
$ws.Range("E47").Value = "  -7.39%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.48"
$ws.Range("E48").Value = "  +4.23%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "10.25"
$ws.Range("E49").Value = "  +7.89%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.16"
$ws.Range("E50").Value = "  -6.04%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.08"
$ws.Range("E51").Value = "  -6.11%  "
